# edit.ps1 -- apply the "added fastText instead of langdetect" edits
#
# Summary of changes (from the OOXML diff):
#   1. Merge the two runs "Pastrand doar adjective si substantive" + ":"
#      into a single run "Pastrand doar adjective si substantive:".
#   2. Remove the standalone paragraph "Id = 562:" that used to sit right
#      before the "829 - DATCU Mihai" paragraph.
#   3. Move the <w:lastRenderedPageBreak/> marker that was on the
#      "841 - ANDRONESCU Ecaterina" paragraph down onto the following
#      "0: 0.005*...x ray diffraction xrd..." paragraph.
#   4. Move the <w:lastRenderedPageBreak/> marker that was on the
#      "69354 - VLAD MAGDALENA" paragraph down onto the following
#      "0: 0.020*...fast ion..." paragraph.

$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $pattern) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

# Builds the opening <w:p ...> attribute list (paraId/textId/rsid*) for a
# paragraph, pulled live off the paragraph's own Range.WordOpenXML, so the
# replacement paragraph keeps its original identity attributes.
function Get-ParaOpenAttrs($doc, $paraIndex) {
    $xml = $doc.Paragraphs($paraIndex).Range.WordOpenXML
    if ($xml -match '<w:p(?:\s[^>]*)?>') {
        $tag = $matches[0]
        $attrs = $tag.Substring(4, $tag.Length - 5).Trim()
        return $attrs
    }
    return ""
}

# Minimal XML-escape for text placed inside a <w:t> element.
function Escape-XmlText($text) {
    $text = $text -replace '&', '&amp;'
    $text = $text -replace '<', '&lt;'
    $text = $text -replace '>', '&gt;'
    return $text
}

# ---------------------------------------------------------------------
# 1. Merge "Pastrand doar adjective si substantive" + ":" into one run.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Pastrand doar adjective si substantive:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Pastrand doar adjective si substantive:", 2)

# ---------------------------------------------------------------------
# 2. Delete the "Id = 562:" paragraph that precedes "829 - DATCU Mihai".
# ---------------------------------------------------------------------
$datcuIdx = Get-ParaIndexByText $d "*829 - DATCU Mihai*"
$idParaIdx = $datcuIdx - 1
if ($d.Paragraphs($idParaIdx).Range.Text -like "*Id = 562:*") {
    $d.Paragraphs($idParaIdx).Range.Delete()
}

# ---------------------------------------------------------------------
# 3 & 4. Move <w:lastRenderedPageBreak/> from the name paragraph onto the
#         following "0: ..." topic-words paragraph.
# ---------------------------------------------------------------------
function Move-PageBreakToNextParagraph($doc, $namePattern) {
    $nameIdx = Get-ParaIndexByText $doc $namePattern
    if ($nameIdx -eq -1) { return }
    $nextIdx = $nameIdx + 1

    $nameText = $doc.Paragraphs($nameIdx).Range.Text.TrimEnd("`r","`n")
    $nameAttrs = Get-ParaOpenAttrs $doc $nameIdx
    $nextText = $doc.Paragraphs($nextIdx).Range.Text.TrimEnd("`r","`n")
    $nextAttrs = Get-ParaOpenAttrs $doc $nextIdx

    $w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

    $nameSpace = ""
    if ($nameText -match "^\s" -or $nameText -match "\s$") { $nameSpace = ' xml:space="preserve"' }
    $nextSpace = ""
    if ($nextText -match "^\s" -or $nextText -match "\s$") { $nextSpace = ' xml:space="preserve"' }

    $nameTextX = Escape-XmlText $nameText
    $nextTextX = Escape-XmlText $nextText

    $nameXml = "<w:p $w $nameAttrs><w:r><w:t$nameSpace>$nameTextX</w:t></w:r></w:p>"
    $doc.Paragraphs($nameIdx).Range.InsertXML($nameXml)

    $nextXml = "<w:p $w $nextAttrs><w:r><w:lastRenderedPageBreak/><w:t$nextSpace>$nextTextX</w:t></w:r></w:p>"
    $doc.Paragraphs($nextIdx).Range.InsertXML($nextXml)
}

Move-PageBreakToNextParagraph $d "*841 - ANDRONESCU Ecaterina*"
Move-PageBreakToNextParagraph $d "*69354 - VLAD MAGDALENA*"

Write-Output "done"
